# Apply updated crypto price/volume (and row shift for new "Solana" entry /
# removal of "OKB" / addition of "WEMIXTOKEN") data scraped on 2023-04-13.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "30.405.38"
$ws.Range("E2").Value = "  +1.22%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.010.02"
$ws.Range("E3").Value = "  +4.88%  "

# Row 4: TetherUSD
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "

# Row 5: BNB
$ws.Range("D5").Value = "'324.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.38%  "

# Row 6: USDC
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "

# Row 7: XRP
$ws.Range("E7").Value = "  +1.25%  "

# Row 8: Cardano
$ws.Range("D8").Value = "'0.4163"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.60%  "

# Row 9: Dogecoin
$ws.Range("D9").Value = "'0.08801"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.76%  "

# Row 10: Polygon
$ws.Range("E10").Value = "  +2.53%  "

# Row 11: Solana
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "'24.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.78%  "

# Row 12: WrappedEther
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "2.002.14"
$ws.Range("E12").Value = "  +4.49%  "

# Row 13: Polkadot
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'6.590"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.71%  "

# Row 14: Chainlink
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'7.483"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.70%  "

# Row 15: BinanceUSD
$ws.Range("B15").Value = "BinanceUSD"
$ws.Range("C15").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D15").Value = "'1.002"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.01%  "

# Row 16: Litecoin
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "'94.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.31%  "

# Row 17: ShibaInu
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.00001118"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.00%  "

# Row 18: TRON
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.06533"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.54%  "

# Row 19: Avalanche
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "'18.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.86%  "

# Row 20: Dai
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.00%  "

# Row 21: Uniswap
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'6.225"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.84%  "

# Row 22: WrappedBTC
$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value = "30.460.70"
$ws.Range("E22").Value = "  +1.26%  "

# Row 23: Cosmos
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "'12.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.16%  "

# Row 24: Toncoin
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "'2.227"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.31%  "

# Row 25: WrappedliquidstakedEther2.0
$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").Value = "2.244.07"
$ws.Range("E25").Value = "  +4.91%  "

# Row 26: EthereumClassic
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'22.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.16%  "

# Row 27: Monero
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'162.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.58%  "

# Row 28: LidoDAOToken
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.422"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.50%  "

# Row 29: BitcoinCash
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "'131.57"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.00%  "

# Row 30: ImmutableX
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'1.137"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.67%  "

# Row 31: Stellar
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.1053"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.48%  "

# Row 32: Filecoin
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'6.126"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.17%  "

# Row 33: HuobiToken
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "'3.831"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.49%  "

# Row 34: ARBITRUM
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.350"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.55%  "

# Row 35: VeChain
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").Value = "'0.02525"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.32%  "

# Row 36: InternetComputer(DFINITY)
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "'5.457"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.02%  "

# Row 37: Hedera
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.06617"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.12%  "

# Row 38: Aptos
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value = "'12.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.28%  "

# Row 39: FraxShare
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'9.125"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.66%  "

# Row 40: Algorand
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "'0.2200"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.79%  "

# Row 41: TheSandbox
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6663"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.48%  "

# Row 42: TrustWalletToken
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.234"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.14%  "

# Row 43: EnergySwap
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'13.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.00%  "

# Row 44: Decentraland
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").Value = "'0.6176"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.54%  "

# Row 45: NEARProtocol
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "'2.197"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.10%  "

# Row 46: PancakeSwap
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "'3.668"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.87%  "

# Row 47: EOS
$ws.Range("B47").Value = "EOS"
$ws.Range("C47").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D47").Value = "'1.268"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.42%  "

# Row 48: Quant
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'124.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.77%  "

# Row 49: Aave
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'81.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.22%  "

# Row 50: Cronos
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.06897"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.54%  "

# Row 51: WEMIXTOKEN
$ws.Range("B51").Value = "WEMIXTOKEN"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").Value = "'1.109"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.54%  "

